# Add season record columns (Wins, Losses, Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching the style/formatting of the existing header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values, same for every player row (2-52).
$ws.Range("AD2:AD52").Value = 88
$ws.Range("AE2:AE52").Value = 74
$ws.Range("AF2:AF52").Value = 0
